$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit is a permutation of the data rows 2-17 (columns D, J, K, L, M, P),
# while all other columns (A, B, C, E, F, G, H, I, N, O, Q, R) and row 18 stay
# as-is. Row `r` ends up with the D/J/K/L/M/P values that originally belonged
# to row `mapping[r]`.
$mapping = @{
    2  = 4
    3  = 2
    4  = 5
    5  = 9
    6  = 7
    7  = 17
    8  = 14
    9  = 8
    10 = 13
    11 = 3
    12 = 15
    13 = 6
    14 = 16
    15 = 10
    16 = 11
    17 = 12
}

# Snapshot original values for columns D, J, K, L, M, P (rows 2-17) before
# overwriting anything, so the permutation reads consistently from the
# original state rather than partially-updated values.
$cols = @("D", "J", "K", "L", "M", "P")
$original = @{}
foreach ($row in 2..17) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $original[$row] = $rowVals
}

foreach ($row in 2..17) {
    $srcRow = $mapping[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $original[$srcRow][$col]
    }
}
